$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 709.125
$ws.Range("I18").Value = 709.125
$ws.Range("K18").Value = 709.125
$ws.Range("M18").Value = -425.125

$ws.Range("H28").Value = 570.2059
$ws.Range("I28").Value = 436.65384
$ws.Range("J28").Value = 1004.25
$ws.Range("K28").Value = 436.65384
$ws.Range("L28").Value = 1004.25
$ws.Range("M28").Value = 48.34616
$ws.Range("N28").Value = -1974.25

$ws.Range("H38").Value = 102.666664
$ws.Range("J38").Value = 100
$ws.Range("L38").Value = 300
$ws.Range("N38").Value = -1044

$ws.Range("H62").Value = 2666.6667
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 2666.6667
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -16240

$ws.Range("H74").Value = 7123.7646
$ws.Range("I74").Value = 7123.7646
$ws.Range("K74").Value = 7123.7646
$ws.Range("M74").Value = -6187.7646

$ws.Range("H77").Value = 7123.7646
$ws.Range("I77").Value = 7123.7646
$ws.Range("K77").Value = 35618.823
$ws.Range("M77").Value = -30938.823

$ws.Range("H98").Value = 921.2273
$ws.Range("I98").Value = 863.35
$ws.Range("K98").Value = 863.35
$ws.Range("M98").Value = 634.65

$ws.Range("H112").Value = 1380.25
$ws.Range("I112").Value = 1647.8
$ws.Range("J112").Value = 1291.0667
$ws.Range("K112").Value = 4943.4
$ws.Range("L112").Value = 3873.2001
$ws.Range("M112").Value = -3835.4
$ws.Range("N112").Value = -6089.2001

$ws.Range("H122").Value = 921.2273
$ws.Range("I122").Value = 863.35
$ws.Range("K122").Value = 2590.05
$ws.Range("M122").Value = -140.0500000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 497.35715
$ws.Range("I5").Value = 85.375
$ws.Range("J5").Value = 1046.6666
$ws.Range("K5").Value = 85.375
$ws.Range("L5").Value = 1046.6666
$ws.Range("M5").Value = 26.625
$ws.Range("N5").Value = -1270.6666

$ws.Range("H61").Value = 14099037
$ws.Range("I61").Value = 5351808
$ws.Range("K61").Value = 5351808
$ws.Range("M61").Value = -5351596

$ws.Range("H134").Value = 74999.5
$ws.Range("J134").Value = 74999.5
$ws.Range("L134").Value = 74999.5
$ws.Range("N134").Value = -85139.5

$ws.Range("H136").Value = 14099037
$ws.Range("I136").Value = 5351808
$ws.Range("K136").Value = 16055424
$ws.Range("M136").Value = -16052874

$ws.Range("H141").Value = 86429
$ws.Range("J141").Value = 86429
$ws.Range("L141").Value = 86429
$ws.Range("N141").Value = -96789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 497.35715
$ws.Range("I4").Value = 85.375
$ws.Range("J4").Value = 1046.6666
$ws.Range("K4").Value = 85.375
$ws.Range("L4").Value = 1046.6666
$ws.Range("M4").Value = 29.625
$ws.Range("N4").Value = -1276.6666

$ws.Range("H36").Value = 1325.75
$ws.Range("I36").Value = 1325.75
$ws.Range("K36").Value = 1325.75
$ws.Range("M36").Value = -791.75

$ws.Range("H105").Value = 5684.8335
$ws.Range("I105").Value = 5684.8335
$ws.Range("K105").Value = 5684.8335
$ws.Range("M105").Value = -3937.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 965
$ws.Range("I22").Value = 513
$ws.Range("K22").Value = 513
$ws.Range("M22").Value = -163

$ws.Range("H31").Value = 482638.9
$ws.Range("I31").Value = 1387694.5
$ws.Range("J31").Value = 4970.6387
$ws.Range("K31").Value = 1387694.5
$ws.Range("L31").Value = 4970.6387
$ws.Range("M31").Value = -1387399.5
$ws.Range("N31").Value = -5560.6387

$ws.Range("H34").Value = 482638.9
$ws.Range("I34").Value = 1387694.5
$ws.Range("J34").Value = 4970.6387
$ws.Range("K34").Value = 1387694.5
$ws.Range("L34").Value = 4970.6387
$ws.Range("M34").Value = -1387492.5
$ws.Range("N34").Value = -5374.6387

$ws.Range("H132").Value = 5697.8335
$ws.Range("I132").Value = 5691.846
$ws.Range("J132").Value = 5713.4
$ws.Range("K132").Value = 17075.538
$ws.Range("L132").Value = 17140.2
$ws.Range("M132").Value = -14545.538
$ws.Range("N132").Value = -22200.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2033980.9
$ws.Range("I5").Value = 1623833.1
$ws.Range("J5").Value = 3323016.5
$ws.Range("K5").Value = 4871499.300000001
$ws.Range("L5").Value = 9969049.5
$ws.Range("M5").Value = -4871387.300000001
$ws.Range("N5").Value = -9969273.5

$ws.Range("H7").Value = 1250
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -1388
$ws.Range("N7").Value = -6224

$ws.Range("H10").Value = 2008
$ws.Range("I10").Value = 13.333333
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 39.999999
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 99.000001
$ws.Range("N10").Value = -15278

$ws.Range("H37").Value = 50000
$ws.Range("J37").Value = 50000
$ws.Range("L37").Value = 150000
$ws.Range("N37").Value = -150224

$ws.Range("H120").Value = 15255.556
$ws.Range("I120").Value = 9060
$ws.Range("K120").Value = 27180
$ws.Range("M120").Value = -22342

$ws.Range("H121").Value = 7157.1055
$ws.Range("J121").Value = 10791.667
$ws.Range("L121").Value = 32375.001
$ws.Range("N121").Value = -34995.001

$ws.Range("H131").Value = 9489.6
$ws.Range("I131").Value = 4949.5
$ws.Range("K131").Value = 14848.5
$ws.Range("M131").Value = -9808.5

$ws.Range("H135").Value = 2033980.9
$ws.Range("I135").Value = 1623833.1
$ws.Range("J135").Value = 3323016.5
$ws.Range("K135").Value = 14614497.9
$ws.Range("L135").Value = 29907148.5
$ws.Range("M135").Value = -14611962.9
$ws.Range("N135").Value = -29912218.5

$ws.Range("H136").Value = 9743.166999999999
$ws.Range("I136").Value = 4486.3335
$ws.Range("K136").Value = 13459.0005
$ws.Range("M136").Value = -8359.000499999998

$ws.Range("H137").Value = 9180.421
$ws.Range("J137").Value = 11178.538
$ws.Range("L137").Value = 33535.614
$ws.Range("N137").Value = -43735.614

$ws.Range("H138").Value = 110250
$ws.Range("J138").Value = 15375
$ws.Range("L138").Value = 46125
$ws.Range("N138").Value = -56405

$ws.Range("H139").Value = 10699
$ws.Range("I139").Value = 7637.2
$ws.Range("J139").Value = 12400
$ws.Range("K139").Value = 22911.6
$ws.Range("L139").Value = 37200
$ws.Range("M139").Value = -17771.6
$ws.Range("N139").Value = -47480

$ws.Range("H140").Value = 1991.4
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 10698.7
$ws.Range("I141").Value = 4397.4
$ws.Range("K141").Value = 13192.2
$ws.Range("M141").Value = -8012.199999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 66744.664
$ws.Range("I2").Value = 90981.45
$ws.Range("J2").Value = 93.5
$ws.Range("K2").Value = 90981.45
$ws.Range("L2").Value = 93.5
$ws.Range("M2").Value = -90868.45
$ws.Range("N2").Value = -319.5

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H113").Value = 3471.4285
$ws.Range("I113").Value = 3471.4285
$ws.Range("K113").Value = 3471.4285
$ws.Range("M113").Value = -1301.4285

$ws.Range("H132").Value = 16037.393
$ws.Range("I132").Value = 14438.143
$ws.Range("K132").Value = 43314.429
$ws.Range("M132").Value = -40784.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 999.5
$ws.Range("I16").Value = 999.5
$ws.Range("K16").Value = 999.5
$ws.Range("M16").Value = -829.5

$ws.Range("H122").Value = 6055
$ws.Range("I122").Value = 6767.364
$ws.Range("K122").Value = 20302.092
$ws.Range("M122").Value = -17852.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 487.45456
$ws.Range("I107").Value = 486
$ws.Range("J107").Value = 496.66666
$ws.Range("K107").Value = 1458
$ws.Range("L107").Value = 1489.99998
$ws.Range("M107").Value = -139.1819
$ws.Range("N107").Value = -10114.625

$ws.Range("H113").Value = 1256.1052
$ws.Range("I113").Value = 769.7273
$ws.Range("J113").Value = 1924.875
$ws.Range("K113").Value = 2309.1819
$ws.Range("L113").Value = 5774.625
$ws.Range("M113").Value = -139.1819
$ws.Range("N113").Value = -10114.625
